$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(9, 8).Value = 3417.3333
$ws.Cells.Item(9, 9).Value = 0
$ws.Cells.Item(9, 11).Value = 0
$ws.Cells.Item(9, 13).ClearContents()

$ws.Cells.Item(15, 8).Value = 1835.1904
$ws.Cells.Item(15, 9).Value = 1835.1904
$ws.Cells.Item(15, 11).Value = 5505.5712
$ws.Cells.Item(15, 13).Value = -5336.5712

$ws.Cells.Item(97, 8).Value = 1140.8572
$ws.Cells.Item(97, 9).Value = 489.66666
$ws.Cells.Item(97, 10).Value = 1629.25
$ws.Cells.Item(97, 11).Value = 1468.99998
$ws.Cells.Item(97, 12).Value = 4887.75
$ws.Cells.Item(97, 13).Value = -972.9999800000001
$ws.Cells.Item(97, 14).Value = -5879.75

$ws.Cells.Item(138, 8).Value = 3074.5
$ws.Cells.Item(138, 9).Value = 1641.6786
$ws.Cells.Item(138, 10).Value = 4746.125
$ws.Cells.Item(138, 11).Value = 4925.0358
$ws.Cells.Item(138, 12).Value = 14238.375
$ws.Cells.Item(138, 13).Value = 214.9642000000003
$ws.Cells.Item(138, 14).Value = -24518.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(92, 8).Value = 41999
$ws.Cells.Item(92, 10).Value = 41999
$ws.Cells.Item(92, 12).Value = 41999
$ws.Cells.Item(92, 14).Value = -46991

$ws.Cells.Item(97, 8).Value = 1010.82355
$ws.Cells.Item(97, 9).Value = 1011.5
$ws.Cells.Item(97, 11).Value = 1011.5
$ws.Cells.Item(97, 13).Value = -515.5

$ws.Cells.Item(137, 8).Value = 69995
$ws.Cells.Item(137, 10).Value = 69995
$ws.Cells.Item(137, 12).Value = 69995
$ws.Cells.Item(137, 14).Value = -80195

$ws.Cells.Item(139, 8).Value = 63330.668
$ws.Cells.Item(139, 10).Value = 63330.668
$ws.Cells.Item(139, 12).Value = 63330.668
$ws.Cells.Item(139, 14).Value = -73610.66800000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 3626.8572
$ws.Cells.Item(94, 9).Value = 3077.6
$ws.Cells.Item(94, 11).Value = 3077.6
$ws.Cells.Item(94, 13).Value = -2626.6

$ws.Cells.Item(99, 8).Value = 1136.1666
$ws.Cells.Item(99, 9).Value = 1203.4
$ws.Cells.Item(99, 10).Value = 800
$ws.Cells.Item(99, 11).Value = 1203.4
$ws.Cells.Item(99, 12).Value = 800
$ws.Cells.Item(99, 13).Value = 294.5999999999999
$ws.Cells.Item(99, 14).Value = -3796

$ws.Cells.Item(134, 8).Value = 1979.8276
$ws.Cells.Item(134, 9).Value = 1095.9584
$ws.Cells.Item(134, 10).Value = 6222.4
$ws.Cells.Item(134, 11).Value = 3287.8752
$ws.Cells.Item(134, 12).Value = 18667.2
$ws.Cells.Item(134, 13).Value = -752.8751999999999
$ws.Cells.Item(134, 14).Value = -23737.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 3241.0833
$ws.Cells.Item(5, 9).Value = 1987
$ws.Cells.Item(5, 10).Value = 3868.125
$ws.Cells.Item(5, 11).Value = 5961
$ws.Cells.Item(5, 12).Value = 11604.375
$ws.Cells.Item(5, 13).Value = -5849
$ws.Cells.Item(5, 14).Value = -11828.375

$ws.Cells.Item(135, 8).Value = 3241.0833
$ws.Cells.Item(135, 9).Value = 1987
$ws.Cells.Item(135, 10).Value = 3868.125
$ws.Cells.Item(135, 11).Value = 17883
$ws.Cells.Item(135, 12).Value = 34813.125
$ws.Cells.Item(135, 13).Value = -15348
$ws.Cells.Item(135, 14).Value = -39883.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 8268.546
$ws.Cells.Item(80, 9).Value = 2334.6667
$ws.Cells.Item(80, 10).Value = 10493.75
$ws.Cells.Item(80, 11).Value = 2334.6667
$ws.Cells.Item(80, 12).Value = 10493.75
$ws.Cells.Item(80, 13).Value = -1336.6667
$ws.Cells.Item(80, 14).Value = -12489.75

$ws.Cells.Item(83, 8).Value = 8268.546
$ws.Cells.Item(83, 9).Value = 2334.6667
$ws.Cells.Item(83, 10).Value = 10493.75
$ws.Cells.Item(83, 11).Value = 11673.3335
$ws.Cells.Item(83, 12).Value = 52468.75
$ws.Cells.Item(83, 13).Value = -6681.333500000001
$ws.Cells.Item(83, 14).Value = -62452.75

$ws.Cells.Item(97, 8).Value = 2182.1
$ws.Cells.Item(97, 9).Value = 1479.125
$ws.Cells.Item(97, 11).Value = 1479.125
$ws.Cells.Item(97, 13).Value = -983.125

$ws.Cells.Item(102, 8).Value = 2901344.2
$ws.Cells.Item(102, 9).Value = 3336186.5
$ws.Cells.Item(102, 11).Value = 3336186.5
$ws.Cells.Item(102, 13).Value = -3334564.5

$ws.Cells.Item(113, 8).Value = 1487.7142
$ws.Cells.Item(113, 9).Value = 1178.7778
$ws.Cells.Item(113, 10).Value = 2043.8
$ws.Cells.Item(113, 11).Value = 1178.7778
$ws.Cells.Item(113, 12).Value = 2043.8
$ws.Cells.Item(113, 13).Value = 991.2221999999999
$ws.Cells.Item(113, 14).Value = -6383.8

$ws.Cells.Item(126, 8).Value = 3360.9473
$ws.Cells.Item(126, 9).Value = 2110.8333
$ws.Cells.Item(126, 11).Value = 6332.499899999999
$ws.Cells.Item(126, 13).Value = -3862.499899999999

$ws.Cells.Item(133, 8).Value = 68795.42999999999
$ws.Cells.Item(133, 10).Value = 68795.42999999999
$ws.Cells.Item(133, 12).Value = 68795.42999999999
$ws.Cells.Item(133, 14).Value = -78915.42999999999

$ws.Cells.Item(140, 8).Value = 68778.664
$ws.Cells.Item(140, 10).Value = 68778.664
$ws.Cells.Item(140, 12).Value = 68778.664
$ws.Cells.Item(140, 14).Value = -79138.664

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1799.8572
$ws.Cells.Item(22, 9).Value = 1000
$ws.Cells.Item(22, 10).Value = 1933.1666
$ws.Cells.Item(22, 11).Value = 1000
$ws.Cells.Item(22, 12).Value = 1933.1666
$ws.Cells.Item(22, 13).Value = -705
$ws.Cells.Item(22, 14).Value = -2523.1666

$ws.Cells.Item(27, 8).Value = 1799.8572
$ws.Cells.Item(27, 9).Value = 1000
$ws.Cells.Item(27, 10).Value = 1933.1666
$ws.Cells.Item(27, 11).Value = 1000
$ws.Cells.Item(27, 12).Value = 1933.1666
$ws.Cells.Item(27, 13).Value = -893
$ws.Cells.Item(27, 14).Value = -2147.1666

$ws.Cells.Item(46, 8).Value = 2699.3
$ws.Cells.Item(46, 9).Value = 998.5
$ws.Cells.Item(46, 11).Value = 998.5
$ws.Cells.Item(46, 13).Value = -810.5

$ws.Cells.Item(61, 8).Value = 5541
$ws.Cells.Item(61, 9).Value = 2720.3333
$ws.Cells.Item(61, 11).Value = 2720.3333
$ws.Cells.Item(61, 13).Value = -2518.3333

$ws.Cells.Item(68, 8).Value = 98229.55
$ws.Cells.Item(68, 9).Value = 254305.5
$ws.Cells.Item(68, 11).Value = 254305.5
$ws.Cells.Item(68, 13).Value = -253556.5

$ws.Cells.Item(71, 8).Value = 98229.55
$ws.Cells.Item(71, 9).Value = 254305.5
$ws.Cells.Item(71, 11).Value = 1271527.5
$ws.Cells.Item(71, 13).Value = -1267783.5

$ws.Cells.Item(94, 8).Value = 49999.5
$ws.Cells.Item(94, 10).Value = 49999.5
$ws.Cells.Item(94, 12).Value = 49999.5
$ws.Cells.Item(94, 14).Value = -51351.5

$ws.Cells.Item(100, 8).Value = 8953.875
$ws.Cells.Item(100, 9).Value = 4682.7856
$ws.Cells.Item(100, 11).Value = 4682.7856
$ws.Cells.Item(100, 13).Value = -4141.7856

$ws.Cells.Item(113, 8).Value = 5541
$ws.Cells.Item(113, 9).Value = 2720.3333
$ws.Cells.Item(113, 11).Value = 2720.3333
$ws.Cells.Item(113, 13).Value = -550.3332999999998

$ws.Cells.Item(122, 8).Value = 13165.333
$ws.Cells.Item(122, 10).Value = 8501.25
$ws.Cells.Item(122, 12).Value = 25503.75
$ws.Cells.Item(122, 14).Value = -30403.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 4388
$ws.Cells.Item(62, 9).Value = 4387
$ws.Cells.Item(62, 11).Value = 4387
$ws.Cells.Item(62, 13).Value = -3763

$ws.Cells.Item(65, 8).Value = 4388
$ws.Cells.Item(65, 9).Value = 4387
$ws.Cells.Item(65, 11).Value = 21935
$ws.Cells.Item(65, 13).Value = -18815

$ws.Cells.Item(70, 8).Value = 38749.75
$ws.Cells.Item(70, 10).Value = 41499.5
$ws.Cells.Item(70, 12).Value = 41499.5
$ws.Cells.Item(70, 14).Value = -42129.5

$ws.Cells.Item(73, 8).Value = 38749.75
$ws.Cells.Item(73, 10).Value = 41499.5
$ws.Cells.Item(73, 12).Value = 41499.5
$ws.Cells.Item(73, 14).Value = -43683.5

$ws.Cells.Item(113, 8).Value = 378.875
$ws.Cells.Item(113, 9).Value = 216.33333
$ws.Cells.Item(113, 10).Value = 866.5
$ws.Cells.Item(113, 11).Value = 648.99999
$ws.Cells.Item(113, 12).Value = 2599.5
$ws.Cells.Item(113, 13).Value = 1521.00001
$ws.Cells.Item(113, 14).Value = -6939.5

$ws.Cells.Item(126, 8).Value = 2414.7
$ws.Cells.Item(126, 9).Value = 1902.6
$ws.Cells.Item(126, 11).Value = 5707.799999999999
$ws.Cells.Item(126, 13).Value = -3237.799999999999
